$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.865.46"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.813.42"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "308.77"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.4661"
$ws.Range("E7").Value = "  +1.90%  "
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").Value = "0.07370"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "0.8699"
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("D11").Value = "20.40"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "1.862.75"
$ws.Range("E12").Value = "  +3.69%  "
$ws.Range("D13").Value = "5.349"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "0.07071"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "91.71"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "6.499"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "0.000008690"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "14.72"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "26.912.28"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "5.331"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("D23").Value = "10.56"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").Value = "2.039.91"
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").Value = "150.89"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").Value = "2.175"
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("D28").Value = "18.34"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("D29").Value = "5.313"
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("D30").Value = "115.65"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("D31").Value = "0.08930"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").Value = "0.7658"
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "4.504"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("D35").Value = "2.902"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "1.088"
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("D38").Value = "0.01960"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "0.05280"
$ws.Range("D40").Value = "2.939"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("D41").Value = "7.272"
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("D43").Value = "2.363"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").Value = "0.1663"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").Value = "8.416"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("D46").Value = "0.4925"
$ws.Range("E46").Value = "  -2.23%  "
$ws.Range("D47").Value = "10.46"
$ws.Range("E47").Value = "  +2.40%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").Value = "103.61"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").Value = "0.06287"
$ws.Range("E51").Value = "  -0.05%  "
